# Weekly crime-stat refresh: roll the report forward one week
# (Volume/Number label + "Week Covering" date range) and replace the
# precinct crime-complaint table (rows 15-28) with the newly collected
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 31   Number  12" -> "...  13" -------------------
$ws.Range("A8").Value = "Volume 31   Number  13"

# --- Header: report week date range -----------------------------------
$ws.Range("C9").Value = "Report Covering the Week  3/25/2024  Through  3/31/2024"

# --- Crime complaints table (rows 15-28) -------------------------------
# Row 15
$ws.Range("L15").Value = -87.5

# Row 16 (C16 was a blank-dash placeholder ("0") -> now an actual count;
# switch it to the same numeric style used by its neighboring cells)
$ws.Range("C16").Value = 2
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -58.333333333333
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = -20.454545454545
$ws.Range("L16").Value = -23.913043478260
$ws.Range("M16").Value = -25.531914893617
$ws.Range("N16").Value = -85.169491525423

# Row 17
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -4.761904761904
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 61
$ws.Range("K17").Value = -24.590163934426
$ws.Range("L17").Value = -14.814814814814
$ws.Range("M17").Value = 76.923076923076
$ws.Range("N17").Value = -20.689655172413

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 52
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -20
$ws.Range("M18").Value = -40.229885057471
$ws.Range("N18").Value = -86.802030456852

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -5.263157894736
$ws.Range("F19").Value = 81
$ws.Range("G19").Value = 88
$ws.Range("H19").Value = -7.954545454545
$ws.Range("I19").Value = 236
$ws.Range("J19").Value = 249
$ws.Range("K19").Value = -5.220883534136
$ws.Range("L19").Value = 8.256880733944
$ws.Range("M19").Value = -31.594202898550
$ws.Range("N19").Value = -62.951334379905

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = -38.461538461538
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = 14.285714285714
$ws.Range("N20").Value = -97.452229299363

# Row 21
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -12.903225806451
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -12.318840579710
$ws.Range("I21").Value = 379
$ws.Range("J21").Value = 421
$ws.Range("K21").Value = -9.976247030878
$ws.Range("L21").Value = -5.955334987593
$ws.Range("M21").Value = -26.264591439688
$ws.Range("N21").Value = -76.960486322188

# Row 22
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -66.666666666666
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = -21.739130434782
$ws.Range("L22").Value = -37.931034482758

# Row 23
$ws.Range("M23").Value = -60

# Row 24
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 30.769230769230
$ws.Range("F24").Value = 251
$ws.Range("G24").Value = 148
$ws.Range("H24").Value = 69.594594594594
$ws.Range("I24").Value = 781
$ws.Range("J24").Value = 473
$ws.Range("K24").Value = 65.116279069767
$ws.Range("L24").Value = 41.229656419529
$ws.Range("M24").Value = 85.510688836104

# Row 25
$ws.Range("C25").Value = 43
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = 26.470588235294
$ws.Range("F25").Value = 220
$ws.Range("G25").Value = 118
$ws.Range("H25").Value = 86.440677966101
$ws.Range("I25").Value = 667
$ws.Range("J25").Value = 357
$ws.Range("K25").Value = 86.834733893557
$ws.Range("L25").Value = 63.882063882063

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 57.142857142857
$ws.Range("F26").Value = 48
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 141
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = 28.181818181818
$ws.Range("L26").Value = 23.684210526315
$ws.Range("M26").Value = 58.426966292134

# Row 27
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -57.142857142857
$ws.Range("L27").Value = -72.727272727272

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 12
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 26
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = -10.344827586206
$ws.Range("L28").Value = 13.043478260869
